$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  ,@(0.209860093271153,0.120586275816123,0.0201532311792139,0.716688874083944,0.163724183877415,0.06545636242505,0.841938707528314,0.00116588940706196,0.963857428381079,0.950532978014657,0.000499666888740839,0.00266489007328448,0.985676215856096,0,0.00916055962691539,0,0.997834776815456,0.117754830113258,0.992838107928048,0.38124583610926,0.889240506329114,0.0564623584277149,0.0389740173217855)
  ,@(0.122251832111925,0.851598934043971,0.0709526982011992,0.044636908727515,0.797301798800799,0.869420386409061,0.00016655562958028,0.0541305796135909,0.00582944703530979,0.00316455696202532,0.0438041305796136,0.0326449033977348,0.00016655562958028,0.00016655562958028,0.000499666888740839,1,0.00016655562958028,0,0.00283144570286476,0.0016655562958028,0.0024983344437042,0.00233177881412392,0.00216522318454364)
  ,@(0.626915389740173,0.0153231179213857,0.00199866755496336,0.209860093271153,0.0146568954030646,0.00832778147901399,0.146568954030646,0,0.0259826782145237,0.0436375749500333,0.000999333777481679,0.957361758827448,0.00516322451698867,0,0.989173884077282,0,0.000999333777481679,0.878247834776815,0.00366422385076616,0.564790139906729,0.101598934043971,0.895736175882745,0.914890073284477)
  ,@(0.0269820119920053,0.00766155896069287,0.898401065956029,0.0144903397734843,0.0174883411059294,0.048800799467022,0.00183211192538308,0.943537641572285,0.000832778147901399,0.00199866755496336,0.954696868754164,0.00582944703530979,0.00849433710859427,0.99983344437042,0.00016655562958028,0,0.000999333777481679,0,0,0.0433044636908727,0.00133244503664224,0.0431379080612925,0.043471019320453)
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = 2 + $i
    $rowData = $data[$i]
    for ($j = 0; $j -lt $rowData.Length; $j++) {
        $col = 2 + $j
        $ws.Cells.Item($row, $col).Value = $rowData[$j]
    }
}

Write-Output "Done updating frequency table values"